# Insert a new review row at row 15 (pushing existing rows 15-24 down to 16-25)
# and populate it with the new review data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("15:15").Insert()

$ws.Range("A15").Value = 5
$ws.Range("B15").Value = "Solução rápida e competente "
$ws.Range("C15").Value = 45958.72571003472
$ws.Range("D15").Value = "OTcxNWNiMjktMWZhMC00NTg2LWExNmYtYWVlYzAwNDkyMWNlOjU3MDE2"
